# Rename the original (only) sheet to "Hoja2" and add three more sheets
# ("Sheet" with the song data, then the empty "Hoja3" and "Hoja1"),
# matching: Hoja2, Sheet, Hoja3, Hoja1 (left to right).

$wb = $excel.ActiveWorkbook

# The workbook starts with a single sheet named "Sheet" -> rename to "Hoja2".
$hoja2 = $wb.Worksheets.Item(1)
$hoja2.Name = "Hoja2"

# Insert the new "Sheet" right after "Hoja2" and fill it with the lyrics table.
$sheet = $wb.Worksheets.Add($null, $hoja2)
$sheet.Name = "Sheet"

$sheet.Range("A1").Value = "Canciones"
$sheet.Range("B1").Value = "Género"
$sheet.Range("C1").Value = "Año"
$sheet.Range("D1").Value = "¿Es famosa la canción?"

$sheet.Range("A2").Value = "Vamos a cantar"
$sheet.Range("B2").Value = "¿Qué canción?"
$sheet.Range("C2").Value = "Pues yo no sé, dime tú"
$sheet.Range("D2").Value = "Ok!... entonces cantaremos los pollitos"

# Add the two remaining (empty) sheets after "Sheet", in order: Hoja3, Hoja1.
$hoja3 = $wb.Worksheets.Add($null, $sheet)
$hoja3.Name = "Hoja3"

$hoja1 = $wb.Worksheets.Add($null, $hoja3)
$hoja1.Name = "Hoja1"
